# Automatische test-sync: 2025-06-29 14:58:50
# Append a new test-mail row to the "Logs" sheet, extend the conditional
# formatting that tracks it, roll the corresponding category tally into the
# "Dashboard" sheet, and extend the dashboard chart's source ranges to match.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 23 -------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A23").Value = "Kun je contact opnemen met de klant?"
$logs.Range("B23").Value = "mailmind.test@zohomail.eu"
$logs.Range("C23").Value = "Testmail #8: Kun je contact opnemen met de klant?"
$logs.Range("D23").Value = "Klacht / Probleem"
$logs.Range("E23").Value = "Beste klantenservice,
Ik heb zojuist een testmail verstuurd (Testmail #14) om te controleren of onze klanten goed bereikbaar zijn via e-mail. Kun je bevestigen of deze testmail succesvol is ontvangen en of jullie contact hebben opgenomen met de klant?
Met vriendelijke groet,
[Naam]"
$logs.Range("F23").Value = "2025-06-29 14:57:51"
$logs.Range("G23").Value = "Ja"
$logs.Range("H23").Value = "Nee"
$logs.Range("I23").Value = "Ja"

# Conditional formatting ranges need to grow from row 22 to row 23 along with
# the data (Excel does not auto-expand these just because a cell below was
# filled in).
$logs.Range("D2:D22").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D23"))
$logs.Range("G2:G22").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G23"))
$logs.Range("H2:H22").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H23"))
$logs.Range("I2:I22").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I23"))

# --- Dashboard sheet: append the new "Klacht / Probleem" tally row 7 ----------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A7").Value = "Klacht / Probleem"
$dash.Range("B7").Value = 1

# Extend the bar chart's category/value source ranges to include row 7.
$chart = $dash.ChartObjects(1).Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$7,'Dashboard'!`$B`$2:`$B`$7,1)"
